$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.822.30"
$ws.Range("E2").Value = "  +7.80%  "
$ws.Range("D3").Value = "1.811.60"
$ws.Range("E3").Value = "  +5.01%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'249.47"
$ws.Range("E5").Value = "  +3.65%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "'0.4947"
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("D8").Value = "'0.2782"
$ws.Range("E8").Value = "  +7.86%  "
$ws.Range("D9").Value = "'0.06398"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("D10").Value = "1.814.22"
$ws.Range("E10").Value = "  +5.19%  "
$ws.Range("D11").Value = "'16.75"
$ws.Range("E11").Value = "  +5.11%  "
$ws.Range("D12").Value = "'0.07130"
$ws.Range("E12").Value = "  +3.31%  "
$ws.Range("D13").Value = "'0.6485"
$ws.Range("E13").Value = "  +6.90%  "
$ws.Range("D14").Value = "'84.07"
$ws.Range("E14").Value = "  +9.36%  "
$ws.Range("D15").Value = "'4.697"
$ws.Range("E15").Value = "  +5.01%  "
$ws.Range("D16").Value = "28.789.97"
$ws.Range("E16").Value = "  +8.43%  "
$ws.Range("D17").Value = "'0.9998"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'0.000007399"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("D19").Value = "'0.9998"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'12.23"
$ws.Range("E20").Value = "  +7.04%  "
$ws.Range("D21").Value = "2.048.51"
$ws.Range("E21").Value = "  +5.10%  "
$ws.Range("D22").Value = "'4.594"
$ws.Range("E22").Value = "  +3.81%  "
$ws.Range("D23").Value = "'8.887"
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("D24").Value = "'5.344"
$ws.Range("E24").Value = "  +5.58%  "
$ws.Range("D25").Value = "'143.03"
$ws.Range("E25").Value = "  +4.25%  "
$ws.Range("D26").Value = "'133.15"
$ws.Range("E26").Value = "  +25.44%  "
$ws.Range("E27").Value = "  +8.96%  "
$ws.Range("D28").Value = "'1.887"
$ws.Range("E28").Value = "  +6.72%  "
$ws.Range("D29").Value = "'1.400"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").Value = "'4.160"
$ws.Range("E30").Value = "  +5.62%  "
$ws.Range("D31").Value = "'0.08355"
$ws.Range("E31").Value = "  +5.20%  "
$ws.Range("D32").Value = "'3.848"
$ws.Range("E32").Value = "  +4.31%  "
$ws.Range("D33").Value = "'0.04938"
$ws.Range("E33").Value = "  +10.36%  "
$ws.Range("D34").Value = "'1.089"
$ws.Range("E34").Value = "  +8.32%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.707"
$ws.Range("E35").Value = "  +4.35%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6771"
$ws.Range("E36").Value = "  +9.14%  "
$ws.Range("D37").Value = "'2.270"
$ws.Range("E37").Value = "  +11.99%  "
$ws.Range("D38").Value = "'2.763"
$ws.Range("E38").Value = "  +13.08%  "
$ws.Range("D39").Value = "'0.9541"
$ws.Range("E39").Value = "  +3.04%  "
$ws.Range("D40").Value = "'6.054"
$ws.Range("E40").Value = "  +7.17%  "
$ws.Range("D41").Value = "'0.01589"
$ws.Range("E41").Value = "  +6.37%  "
$ws.Range("D42").Value = "'0.9996"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "'100.64"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").Value = "'0.4094"
$ws.Range("E44").Value = "  +6.87%  "
$ws.Range("D45").Value = "'7.198"
$ws.Range("E45").Value = "  +5.21%  "
$ws.Range("D46").Value = "'0.1222"
$ws.Range("E46").Value = "  +5.78%  "
$ws.Range("D47").Value = "'0.05517"
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("D48").Value = "'8.159"
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("D49").Value = "'31.63"
$ws.Range("E49").Value = "  +5.18%  "
$ws.Range("D50").Value = "'1.311"
$ws.Range("E50").Value = "  +6.67%  "
$ws.Range("D51").Value = "'0.3624"
$ws.Range("E51").Value = "  +8.00%  "
